$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 139, shifting existing rows 139-223 down to 140-224
$ws.Rows.Item(139).Insert()

# Populate the newly inserted row 139 with the new weekly price record
$ws.Range("A139").Value = 9
$ws.Range("B139").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C139").Value = "Metropolitana"
$ws.Range("D139").Value = 44596
$ws.Range("D139").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E139").Value = 13
$ws.Range("F139").Value = 100112030
$ws.Range("G139").Value = "Poroto granado"
$ws.Range("H139").Value = "Sin especificar"
$ws.Range("I139").Value = "Primera"
$ws.Range("J139").Value = 160
$ws.Range("K139").Value = 26000
$ws.Range("L139").Value = 27000
$ws.Range("M139").Value = 26562
$ws.Range("N139").Value = "$/saco 25 kilos"
$ws.Range("O139").Value = "Región Metropolitana"
$ws.Range("P139").Value = 1062
$ws.Range("Q139").Value = 25
$ws.Range("R139").Value = "Hortaliza"
